$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7419204711914062
$ws.Range("B1").Value = 1.856776118278503
$ws.Range("C1").Value = 4.718027591705322
$ws.Range("D1").Value = 1.915159583091736
$ws.Range("E1").Value = 1.276799082756042
